$d = $word.ActiveDocument

# 1. "from 2017 and 2019" -> "from 2017 to 2019"
$d.Content.Find.Execute("from 2017 and 2019", $false, $false, $false, $false, $false,
                         $true, 1, $false, "from 2017 to 2019", 2)

# 2. Insert the new sentences. In the source edit, the author's cursor (the
#    _GoBack bookmark) sat right after "...from 2017 to 2019. " and the new
#    text was typed in two bursts, leaving the bookmark in the middle, between
#    "...individual growth and condition " and "were not significantly...".
$beforeText = "Capture-mark-recapture analyses using Barker’s joint live-recapture/tag-recovery model showed that ___________ effect of harvest on natural mortality, while changes to individual growth and condition "
$afterText = "were not significantly related to the degree of biomass removed in all lakes. Commercial harvest of Bigmouth Buffalo from shallow, natural lakes of Iowa should be monitored to ensure this native species is not extirpated within its historic range."

$bmRange = $d.Bookmarks.Item("_GoBack").Range
$bmStart = $bmRange.Start

# Insert the whole new block in one go via InsertBefore so that it correctly
# inherits the surrounding Times New Roman / 12pt run formatting.
$bmRange.InsertBefore($beforeText + $afterText)

# Re-plant the _GoBack bookmark between $beforeText and $afterText, exactly
# where the original edit left the cursor.
$d.Bookmarks.Item("_GoBack").Delete()
$splitPoint = $bmStart + $beforeText.Length
$newBmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $newBmRange)
